# Update "想去人数" (want-to-go count) values in column F across sheets,
# reflecting freshly scraped figures from the site rebuild.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 300
$ws1.Range("F10").Value = 919
$ws1.Range("F18").Value = 357
$ws1.Range("F20").Value = 1330
$ws1.Range("F28").Value = 3373
$ws1.Range("F31").Value = 1497

# Sheet "演出" (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F7").Value = 3

# Sheet "本地生活" (Local Life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 780

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 780
$ws4.Range("F9").Value = 300
$ws4.Range("F16").Value = 3
$ws4.Range("F21").Value = 919
$ws4.Range("F29").Value = 357
$ws4.Range("F31").Value = 1330
$ws4.Range("F41").Value = 3373
$ws4.Range("F44").Value = 1497
